{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// --- Step 1: last three rows (originally multi-value rows) collapse to a single value each ---\ntable.getCell(43, 0).value = \"99.99\";\ntable.getCell(44, 0).value = \"0.08\";\ntable.getCell(45, 0).value = \"865\";\n\n// --- Step 2: rows 10 & 11 (\"0.00005\" and \"0.00014\") merge into a single row \"0.07665\" ---\ntable.getCell(10, 0).value = \"0.07665\";\ntable.deleteRows(11, 1);\n\n// --- Step 3: rows 7,8,9 (\"0.00001\",\"0.00004\",\"0.00005\") are removed entirely ---\ntable.deleteRows(7, 3);\n\n// --- Step 4: row 4 (\"0.00004\") becomes \"0.00002\", and three new rows are inserted after it ---\ntable.getCell(4, 0).value = \"0.00002\";\n\n// --- Step 5: the first four rows get simple value replacements ---\ntable.getCell(0, 0).value = \"0M\";\ntable.getCell(1, 0).value = \"0M\";\ntable.getCell(2, 0).value = \"0M\";\ntable.getCell(3, 0).value = \"1797\";\n\ntable.rows.load(\"items\");\nawait context.sync();\nconst row4 = table.rows.items[4];\nrow4.insertRows(\"After\", 3, [[\"0.00010\"], [\"0.00005\"], [\"0.00001\"]]);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# --- Step 1: last three rows (originally multi-value rows) collapse to a single value each ---\n$t.Cell(44, 1).Range.Text = \"99.99\"\n$t.Cell(45, 1).Range.Text = \"0.08\"\n$t.Cell(46, 1).Range.Text = \"865\"\n\n# --- Step 2: rows 11 & 12 (\"0.00005\" and \"0.00014\") merge into a single row \"0.07665\" ---\n$t.Cell(11, 1).Range.Text = \"0.07665\"\n$t.Rows.Item(12).Delete()\n\n# --- Step 3: rows 8,9,10 (\"0.00001\",\"0.00004\",\"0.00005\") are removed entirely ---\n$t.Rows.Item(10).Delete()\n$t.Rows.Item(9).Delete()\n$t.Rows.Item(8).Delete()\n\n# --- Step 4: row 5 (\"0.00004\") becomes \"0.00002\", and three new rows are inserted after it ---\n$t.Cell(5, 1).Range.Text = \"0.00002\"\n$refRow = $t.Rows.Item(6)\n$newRow = $t.Rows.Add($refRow)\n$newRow.Cells.Item(1).Range.Text = \"0.00001\"\n$newRow = $t.Rows.Add($refRow)\n$newRow.Cells.Item(1).Range.Text = \"0.00005\"\n$newRow = $t.Rows.Add($refRow)\n$newRow.Cells.Item(1).Range.Text = \"0.00010\"\n\n# --- Step 5: the first four rows get simple value replacements ---\n$t.Cell(1, 1).Range.Text = \"0M\"\n$t.Cell(2, 1).Range.Text = \"0M\"\n$t.Cell(3, 1).Range.Text = \"0M\"\n$t.Cell(4, 1).Range.Text = \"1797\"\n"}
